$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (report volume/number + week-covering dates).
# These are rich-text shared strings, but only the textual content matters;
# assigning .Value keeps the owning cell's own style (merged header cells).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# ---------------------------------------------------------------------------
# Helper functions to change a cell's underlying type (number <-> text)
# while keeping it visually styled like its neighbours which already use
# that representation (borrows number format / style via PasteSpecial).
# ---------------------------------------------------------------------------
function Set-TextCell($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("A14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

function Set-NumberCellStyle15($addr, $num) {
    $ws.Range($addr).Value = $num
    $ws.Range("I14").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Cells that flip between numeric and text ("0" / "***.*" placeholders).
# ---------------------------------------------------------------------------
Set-TextCell "C16" "0"
Set-TextCell "D20" "0"
Set-TextCell "E20" "***.*"
Set-TextCell "C23" "0"
Set-TextCell "G30" "0"
Set-TextCell "H30" "***.*"

Set-NumberCellStyle15 "C26" 1
Set-NumberCellStyle15 "F26" 1

# ---------------------------------------------------------------------------
# Plain numeric value updates (new weekly crime-stat figures). Each of these
# keeps its existing cell style / number format untouched.
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 25
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = 21.739130434782
$ws.Range("L16").Value = 64.705882352941
$ws.Range("N16").Value = -81.081081081081

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 48
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 33.333333333333
$ws.Range("L17").Value = 45.454545454545
$ws.Range("M17").Value = 11.627906976744
$ws.Range("N17").Value = -63.076923076923

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -54.545454545454
$ws.Range("I18").Value = 37
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = -38.333333333333
$ws.Range("L18").Value = 85
$ws.Range("M18").Value = -19.565217391304
$ws.Range("N18").Value = -79.891304347826

$ws.Range("C19").Value = 3
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 14
$ws.Range("H19").Value = -6.666666666666
$ws.Range("I19").Value = 68
$ws.Range("J19").Value = 77
$ws.Range("K19").Value = -11.688311688311
$ws.Range("L19").Value = 30.769230769230
$ws.Range("M19").Value = -21.839080459770
$ws.Range("N19").Value = -20.930232558139

$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 32
$ws.Range("K20").Value = 39.130434782608
$ws.Range("L20").Value = 88.235294117647
$ws.Range("M20").Value = 10.344827586206
$ws.Range("N20").Value = -82.608695652173

$ws.Range("C21").Value = 9
$ws.Range("E21").Value = -10
$ws.Range("F21").Value = 38
$ws.Range("G21").Value = 42
$ws.Range("H21").Value = -9.523809523809
$ws.Range("I21").Value = 215
$ws.Range("J21").Value = 221
$ws.Range("K21").Value = -2.714932126696
$ws.Range("L21").Value = 51.408450704225
$ws.Range("M21").Value = -10.788381742738
$ws.Range("N21").Value = -70.985155195681

$ws.Range("M22").Value = -80

$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -25
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 8.888888888888
$ws.Range("L23").Value = 22.5
$ws.Range("M23").Value = 68.965517241379

$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 90.909090909090
$ws.Range("F24").Value = 53
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = 47.222222222222
$ws.Range("I24").Value = 277
$ws.Range("J24").Value = 171
$ws.Range("K24").Value = 61.988304093567
$ws.Range("L24").Value = 88.435374149659
$ws.Range("M24").Value = 41.326530612244

$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -88.888888888888
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -55
$ws.Range("I25").Value = 69
$ws.Range("J25").Value = 84
$ws.Range("K25").Value = -17.857142857142
$ws.Range("L25").Value = 15
$ws.Range("M25").Value = -43.902439024390

$ws.Range("I26").Value = 3
$ws.Range("K26").Value = 50
$ws.Range("L26").Value = -25

$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -46.153846153846

$ws.Range("M28").Value = -50

$ws.Range("M29").Value = -50
